$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Read the existing "File Name" / "Unnormalized P_max" data (rows 2..lastRow)
# ---------------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$rows = @()
for ($i = 2; $i -le $lastRow; $i++) {
    $name = $ws.Cells.Item($i, 1).Value2
    $val  = $ws.Cells.Item($i, 2).Value2

    # Electrode location = the File Name up to the first underscore,
    # e.g. "A11_bipolar_10V_1kHz.txt" -> "A11"
    if ($name -match "^([A-Za-z]+)(\d+)_") {
        $letter = $matches[1]
        $num    = [int]$matches[2]
        $loc    = "$letter$num"
    } else {
        $letter = ""
        $num    = 0
        $loc    = ""
    }

    $rows += [PSCustomObject]@{
        FileName = $name
        Value    = $val
        Letter   = $letter
        Num      = $num
        Loc      = $loc
    }
}

# ---------------------------------------------------------------------------
# Sort by electrode location, A1 -> O15 (letter, then numeric position)
# ---------------------------------------------------------------------------
$sortedRows = $rows | Sort-Object Letter, Num

# ---------------------------------------------------------------------------
# Add the "Electrode Locations" header in column C, matching the style
# already used by the other header cells (A1 / B1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Electrode Locations"

# ---------------------------------------------------------------------------
# Write the sorted rows back out across columns A, B and C
# ---------------------------------------------------------------------------
$r = 2
foreach ($row in $sortedRows) {
    $ws.Cells.Item($r, 1).Value = $row.FileName
    $ws.Cells.Item($r, 2).Value = $row.Value
    $ws.Cells.Item($r, 3).Value = $row.Loc
    $r = $r + 1
}
